{"js": "// Insert three short phrases into two instructional paragraphs:\n//   1. \"...nh\u01b0 /etc/passwd\" -> \"...nh\u01b0 /etc/passwd trong m\u00e1y agent\"  (before the comma)\n//   2. \"...ossec.conf\" -> \"...ossec.conf trong m\u00e1y agent\"            (before \" v\u00e0 th\u00eam\")\n//   3. \"...OSSEC agent:\" -> \"...OSSEC agent v\u00e0 OSSEC server:\"        (before the colon)\n//\n// Each insertion is located with a short, unambiguous search string so the\n// new text lands exactly at the point the edit requires, regardless of how\n// the surrounding sentence is currently split into runs.\n\nasync function insertAfterText(needle, insertion) {\n  const hits = context.document.body.search(needle, { matchCase: true, matchWholeWord: false });\n  hits.load(\"items\");\n  await context.sync();\n  if (hits.items.length === 0) {\n    throw new Error(\"Could not find target text: \" + needle);\n  }\n  // Collapse to a zero-width point immediately after the matched text, then\n  // insert the new phrase right there (keeps the rest of the sentence intact).\n  const insertionPoint = hits.items[0].getRange(\"After\");\n  insertionPoint.insertText(insertion, \"Replace\");\n  await context.sync();\n}\n\n// 1) \"nh\u01b0 /etc/passwd\" is unique in the document, so anchor on that and\n//    insert right after \"/etc/passwd\" (before the following comma).\nawait insertAfterText(\"nh\u01b0 /etc/passwd\", \" trong m\u00e1y agent\");\n\n// 2) \"/var/ossec/etc/ossec.conf\" is unique in the document.\nawait insertAfterText(\"/var/ossec/etc/ossec.conf\", \" trong m\u00e1y agent\");\n\n// 3) \"Kh\u1edfi \u0111\u1ed9ng l\u1ea1i OSSEC agent\" is unique in the document; insert right\n//    after it (before the following colon).\nawait insertAfterText(\"Kh\u1edfi \u0111\u1ed9ng l\u1ea1i OSSEC agent\", \" v\u00e0 OSSEC server\");\n", "ps1": "# Insert three short phrases into two instructional paragraphs:\n#   1. \"...nh\u01b0 /etc/passwd\" -> \"...nh\u01b0 /etc/passwd trong m\u00e1y agent\"  (before the comma)\n#   2. \"...ossec.conf\" -> \"...ossec.conf trong m\u00e1y agent\"            (before \" v\u00e0 th\u00eam\")\n#   3. \"...OSSEC agent:\" -> \"...OSSEC agent v\u00e0 OSSEC server:\"        (before the colon)\n#\n# Each insertion is located with a short, unambiguous Find.Execute target so\n# the new text lands exactly at the point the edit requires.\n\n$d = $word.ActiveDocument\n\nfunction Insert-After-Text($doc, $searchText, $insertText) {\n    $rng = $doc.Content\n    $found = $rng.Find.Execute($searchText)\n    if (-not $found) {\n        throw \"Could not find target text: $searchText\"\n    }\n    # Collapse the found range to its end point, then insert right there so\n    # the rest of the sentence stays intact.\n    $rng.Collapse(0)   # wdCollapseEnd\n    $rng.InsertAfter($insertText)\n}\n\n# 1) \"nh\u01b0 /etc/passwd\" is unique in the document; insert right after it\n#    (i.e. right before the following comma).\nInsert-After-Text $d \"nh\u01b0 /etc/passwd\" \" trong m\u00e1y agent\"\n\n# 2) \"/var/ossec/etc/ossec.conf\" is unique in the document.\nInsert-After-Text $d \"/var/ossec/etc/ossec.conf\" \" trong m\u00e1y agent\"\n\n# 3) \"Kh\u1edfi \u0111\u1ed9ng l\u1ea1i OSSEC agent\" is unique in the document; insert right\n#    after it (i.e. right before the following colon).\nInsert-After-Text $d \"Kh\u1edfi \u0111\u1ed9ng l\u1ea1i OSSEC agent\" \" v\u00e0 OSSEC server\"\n"}
